# Fixing "semaine" and "emploi" excel - update class list (ImporterClasse)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F = "Nom" (class name), Column G = "Cycle" (CI/CP)
# Rebuild the full F5:G19 table with the new, corrected list of classes.
$rows = @(
  @("Nom",     "Cycle"),
  @("3.GINFO", "CI"),
  @("4.GINFO", "CI"),
  @("5.GINFO", "CI"),
  @("Cp 1",    "CP"),
  @("Cp 2",    "CP"),
  @("3.GTR",   "CI"),
  @("4.GTR",   "CI"),
  @("5.GTR",   "CI"),
  @("3.GPMC",  "CI"),
  @("4.GPMC",  "CI"),
  @("5.GPMC",  "CI"),
  @("3.INDUS", "CI"),
  @("4.INDUS", "CI"),
  @("5.INDUS", "CI")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = 5 + $i
  $pair = $rows[$i]
  $ws.Cells.Item($r, 6).Value = $pair[0]
  $ws.Cells.Item($r, 7).Value = $pair[1]
}

# Match the active selection saved with the workbook
$ws.Range("F12").Select()

# Best-effort: restore the normal (non-maximized) window position/size
# that Excel records in the workbook view when the file is saved normally.
try {
  $win = $excel.ActiveWindow
  $win.WindowState = -4143  # xlNormal
  $win.Left = 0
  $win.Top = 0
  $win.Width = 20490
  $win.Height = 7905
} catch {
}
